$d = $word.ActiveDocument

# 1. "I" + "ntroduction" -> single run "Introduction"
#    (MatchCase so the lowercase "introduction" elsewhere is left alone)
$find1 = $d.Content.Find
$find1.Execute("Introduction", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction", 2) | Out-Null

# 2. "Quality of wor" + "k" -> single run "Quality of work"
$find2 = $d.Content.Find
$find2.Execute("Quality of work", $true, $false, $false, $false, $false, $true, 1, $false, "Quality of work", 2) | Out-Null

# 3. Append " (storyboard)" right after "What does this mean???" in red
#    Times New Roman text, without the yellow highlight that the
#    preceding run has.
$rng = $d.Content
$find3 = $rng.Find
$find3.Execute("What does this mean???", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" (storyboard)")
$rng.Font.Name = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"
$rng.Font.Color = 255
